$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rkap")
$ws.Range("A1").Value = "test"
